# Fruta / hortaliza, semanal
# This workbook's data rows (2-20) get re-ordered: the values in the
# "varying" columns (D, L, M, N, O, P, S) of each row are replaced by the
# values that (in the original file) belonged to a different row, per the
# mapping below (new row -> source row in the original workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row number -> source row number (as in the ORIGINAL file)
$rowMap = @{
    2  = 12
    3  = 13
    4  = 17
    5  = 18
    6  = 15
    7  = 16
    8  = 14
    9  = 4
    10 = 3
    11 = 10
    12 = 11
    13 = 6
    14 = 7
    15 = 19
    16 = 20
    17 = 8
    18 = 9
    19 = 2
    20 = 5
}

$cols = @("D", "L", "M", "N", "O", "P", "S")

# Snapshot all original values for the columns that vary by row, before
# writing anything, since several rows are sources for other rows.
$original = @{}
for ($r = 2; $r -le 20; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $original[$r] = $rowVals
}

# Now write the new values for every destination row, taken from the
# snapshot of its mapped source row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $original[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcVals[$col]
    }
}
